$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "StatQuery" text (column C) shared by the Cases/Samples/Files rows.
$newStatQuery = @'
MATCH (p:program)<--(s:study)<-[*]-(c:case)<--(demo:demographic)
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (diag:diagnosis)-->(c)
OPTIONAL MATCH (f:file)-[*]->(c)
OPTIONAL MATCH (sf:file)-->(s)
WITH DISTINCT f, sf, samp AS samp, c, demo, diag, s, p
WHERE demo.breed IN ['Cavalier King Charles Spaniel']
RETURN  
    count(distinct p) AS Programs,
    count(distinct s) AS Studies,
    count(distinct c) AS Cases,
    count(distinct samp) AS Samples,
    count(distinct f) AS `Case Files`,
    count(distinct sf) AS `Study Files`
'@
# Here-string keeps a trailing newline; drop it so the text matches exactly.
$newStatQuery = $newStatQuery.TrimEnd("`r`n")

$ws.Range("C2").Value = $newStatQuery
$ws.Range("C3").Value = $newStatQuery
$ws.Range("C4").Value = $newStatQuery

# Column widths: column C narrower + bestFit-sized, column D narrower + bestFit-sized.
$ws.Columns(3).ColumnWidth = 60.5
$ws.Columns(4).ColumnWidth = 49.33333333333333

# View: zoom to 100%, and move the active selection to B4.
$excel.ActiveWindow.Zoom = 100
$ws.Range("B4").Select()
